# Updates cryptos list data (price + 1h volume%) per the Nov 11 2024 GitHub Actions refresh.
# D/E columns hold plain text (prices use "." as thousands AND decimal separator,
# so some look numeric to Excel -- force NumberFormat "@" first so COM keeps the
# literal text, e.g. "1.00", instead of silently coercing to the number 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "87.869.78"
$ws.Range("E2").Value = "  +10.40%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.332.54"
$ws.Range("E3").Value = "  +5.82%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.31%  "

# Row 5: Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.59"
$ws.Range("E5").Value = "  +5.75%  "

# Row 6: BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "637.07"
$ws.Range("E6").Value = "  +2.06%  "

# Row 7: Dogecoin
$ws.Range("E7").Value = "  +22.98%  "

# Row 8: USDC
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.10%  "

# Row 9: XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").Value = "  +5.45%  "

# Row 10: LidoStakedEther
$ws.Range("D10").Value = "3.338.54"
$ws.Range("E10").Value = "  +6.16%  "

# Row 11: Cardano
$ws.Range("E11").Value = "  +5.66%  "

# Row 12: ShibaInu
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000277"
$ws.Range("E12").Value = "  +10.23%  "

# Row 13: TRON
$ws.Range("E13").Value = "  +2.22%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.953.93"
$ws.Range("E14").Value = "  +6.41%  "

# Row 15: Avalanche
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.32"
$ws.Range("E15").Value = "  +9.53%  "

# Row 16: Toncoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  +3.83%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "87.558.84"
$ws.Range("E17").Value = "  +10.31%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "3.337.34"
$ws.Range("E18").Value = "  +6.40%  "

# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.55"
$ws.Range("E19").Value = "  +2.78%  "

# Row 20: SuiNetwork
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.20"
$ws.Range("E20").Value = "  +8.31%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "448.75"
$ws.Range("E21").Value = "  +3.84%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.11"
$ws.Range("E22").Value = "  +0.90%  "

# Row 23: Polkadot
$ws.Range("E23").Value = "  +3.58%  "

# Row 24: LEO
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.41"
$ws.Range("E24").Value = "  +7.92%  "

# Row 25: NEARProtocol
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.37"
$ws.Range("E25").Value = "  +15.64%  "

# Row 26: Aptos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.31"
$ws.Range("E26").Value = "  +14.99%  "

# Row 27: WrappedeETH
$ws.Range("D27").Value = "3.517.50"
$ws.Range("E27").Value = "  +6.58%  "

# Row 28: Litecoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "78.77"
$ws.Range("E28").Value = "  +4.20%  "

# Row 29: PEPE
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000130"
$ws.Range("E29").Value = "  +7.67%  "

# Row 30: Dai
$ws.Range("E30").Value = "  +0.11%  "

# Row 31: Cronos
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.187"
$ws.Range("E31").Value = "  +54.20%  "

# Row 32: Bittensor
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "602.13"
$ws.Range("E32").Value = "  +9.61%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.31"
$ws.Range("E33").Value = "  +5.09%  "

# Row 34: Binance-PegBSC-USD
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.41%  "

# Row 35: Fetch.AI
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +6.11%  "

# Row 36: PancakeSwap
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.05"
$ws.Range("E36").Value = "  +3.74%  "

# Row 37: Kaspa
$ws.Range("E37").Value = "  +2.60%  "

# Row 38: EthereumClassic
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.44"
$ws.Range("E38").Value = "  +2.63%  "

# Row 39: RenderToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.72"
$ws.Range("E39").Value = "  +21.18%  "

# Row 40/41 swapped coins: was PolygonEcosystemToken, now FirstDigitalUSD
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.36%  "

# Row 40/41 swapped coins: was FirstDigitalUSD, now PolygonEcosystemToken
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.419"
$ws.Range("E41").Value = "  +4.87%  "

# Row 42: WhiteBITCoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.38"
$ws.Range("E42").Value = "  +3.11%  "

# Row 43: Stacks
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.06"
$ws.Range("E43").Value = "  +15.02%  "

# Row 44: dogwifhat
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.06"
$ws.Range("E44").Value = "  +14.38%  "

# Row 45: Monero
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "158.41"
$ws.Range("E45").Value = "  -2.87%  "

# Row 46: USDe
$ws.Range("E46").Value = "  +0.02%  "

# Row 47: Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "189.53"
$ws.Range("E47").Value = "  +1.25%  "

# Row 48: OKB
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.13"
$ws.Range("E48").Value = "  +8.75%  "

# Row 49: ImmutableX
$ws.Range("E49").Value = "  +7.39%  "

# Row 50: Mantle
$ws.Range("E50").Value = "  +1.53%  "

# Row 51: InjectiveProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.50"
$ws.Range("E51").Value = "  +9.06%  "
